# Apply the "Auto Update Data" edit:
# 1. Update the "updated date" timestamp in A1.
# 2. Shift the skill-task table (columns A/B/C, rows 201-366) up by one row:
#    each row r (201..366) takes on the values that row r+1 previously had.
#    Row 369 (the last data row) is left untouched since there is no row 370
#    to pull data from.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the header timestamp
$ws.Range("A1").Value = "更新日期：2025.04.05 17:12:18"

# 2. Shift rows 201 -> 366 up by one (pulling values from the row below),
#    processing from the top down so each source row is read before it is
#    itself overwritten.
for ($r = 201; $r -le 366; $r++) {
    $nr = $r + 1
    $a = $ws.Cells.Item($nr, 1).Text
    $b = $ws.Cells.Item($nr, 2).Text
    $c = $ws.Cells.Item($nr, 3).Text
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
}
